$wb = $excel.ActiveWorkbook

$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# Reorder the sheets so that "review_info" comes first and "hotel_info" second
$wsReview.Move($wsHotel)

# References can go stale after a Move(), so re-fetch the hotel_info sheet by name
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Insert a new "State" column into hotel_info right after "Hotel_Name" (before "City")
$wsHotel.Columns("C").Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"
